$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Copy the "Автоматизировано" marker ("V", styled like the other D-column
# checkmarks) into the rows for test cases 1.4.2 and 1.4.3.
$ws.Range("D29").Copy($ws.Range("D30"))
$ws.Range("D29").Copy($ws.Range("D31"))

# Move selection to reflect the last active cell after the edit
$ws.Range("D34").Select()
